$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rename the sheet (defined names referencing it update automatically)
$ws.Name = "What If"

# Set the sheet tab color to purple (RGB 112,48,160 -> FF7030A0)
$ws.Tab.Color = 10498160

# Zoom out to 80%
$excel.ActiveWindow.Zoom = 80

# Move the active selection
$ws.Range("M19").Select() | Out-Null
